# Initial Altitude_adjust P block. Only PV and SP extracted.
# The main control algorithm is not yet implemented.
#
# Adds a new "Orange RX" section (mirroring the existing pin-mapping
# sections on the sheet) describing the PPM signal / +5V / GND wiring.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New section header + rows, following the same B/E/F layout used by the
# other peripheral blocks already on the sheet (e.g. ESCs at row 41).
$ws.Range("B47").Value = "Orange RX"

$ws.Range("E48").Value = " + "
$ws.Range("E49").Value = " - "
$ws.Range("F48").Value = "5V"
$ws.Range("F47").Value = "PD6"
$ws.Range("E47").Value = "Signal (PPM)"

$ws.Range("F49").Value = "GND"

# Match the view state captured in the saved workbook.
$ws.Range("F34").Select()
$excel.ActiveWindow.ScrollRow = 28
